$wb = $excel.ActiveWorkbook

function Update-LangSheet($sheetName, $mdUrl, $xlfUrl, $handbackDateTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
    $ws.Range("C2").Value = "Handed back: in sync with en-US"

    # Newly populated "Latest Target File" (F2) and "Latest Handback File" (G2) columns,
    # each carrying a hyperlink that mirrors the existing handoff-file links.
    $mdName = "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md"
    $xlfName = "e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.$sheetName.xlf"

    $ws.Hyperlinks.Add($ws.Range("F2"), $mdUrl, "", "", $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrl, "", "", $xlfName) | Out-Null

    # Latest Handback DateTime now reflects the actual handback timestamp.
    $ws.Range("H2").Value = $handbackDateTime
}

Update-LangSheet `
    "zh-cn" `
    "https://github.com/OpenLocalizationTest/oltest/blob/2b7338fff35708235e71d882e94de97c35b7ecac/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/61068100a525372f0aeb6d0e1c3fd8988dd6526d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.zh-cn.xlf" `
    "2016-03-13 10:40:24"

Update-LangSheet `
    "de-de" `
    "https://github.com/OpenLocalizationTest/oltest/blob/2b7338fff35708235e71d882e94de97c35b7ecac/e2e/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.md" `
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6c0e2050953e0799cc992ba2a50e9ddd1ffea27f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e1c6acdf-a28e-493a-8d5c-faf9ad2797fe.1715f9dd4a5f1f9b786b4b69e7ed388343a059f8.de-de.xlf" `
    "2016-03-13 10:40:30"
